$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the very end of the document (after the
#    last run of the final paragraph) to the very start of the document
#    (right after the first paragraph's pPr, before the "Simbolismo" run).
# ---------------------------------------------------------------------------

# Drop the existing bookmark wherever it currently lives.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Word's Bookmarks.Add auto-expands a fully collapsed range that sits at the
# very start of the document (offset 0) to cover the first word. To land a
# truly collapsed bookmark at offset 0 we temporarily insert a placeholder
# character, anchor the (now non-zero-offset) bookmark right after it, and
# then delete the placeholder - the bookmark collapses back to offset 0
# without ever passing through the "add at 0" special case.
$placeholder = $d.Range(0, 0)
$placeholder.InsertBefore("X")

$target = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $target)

$cleanup = $d.Range(0, 1)
$cleanup.Delete()

# ---------------------------------------------------------------------------
# 2) Tighten the page margins to 0.5" (36pt / 720 twips) on every side.
# ---------------------------------------------------------------------------

$d.PageSetup.TopMargin = 36
$d.PageSetup.BottomMargin = 36
$d.PageSetup.LeftMargin = 36
$d.PageSetup.RightMargin = 36
